$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New holiday schedule (date serial values, reason shared-string index effectively given via text)
$ws.Range("A3").Value = 45317
$ws.Range("B3").Value = "Republic day"

$ws.Range("A4").Value = 45297
$ws.Range("B4").Value = "Holiday"

$ws.Range("A5").Value = 45304
$ws.Range("B5").Value = "Holiday"

$ws.Range("A6").Value = 45311
$ws.Range("B6").Value = "Holiday"

$ws.Range("A7").Value = 45318
$ws.Range("B7").Value = "Holiday"

$ws.Range("A8").Value = 45306
$ws.Range("B8").Value = "Pongal"

$ws.Range("A3:XFD3").Select()
